$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.693.45"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.885.99"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'247.84"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "'0.4736"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'0.2925"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.06527"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'0.07791"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "1.890.79"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "'0.7365"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'5.248"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "'282.79"
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").Value = "30.791.26"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "'0.000007550"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'0.9995"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "2.140.77"
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'6.263"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'9.225"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").Value = "'163.99"
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -1.19%  "
$ws.Range("D29").Value = "'1.340"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").Value = "'0.09704"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").Value = "'1.493"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D33").Value = "'4.197"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "'0.04844"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("D35").Value = "'1.125"
$ws.Range("D36").Value = "'0.6980"
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "'2.719"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.01907"
$ws.Range("E38").Value = "  +2.09%  "
$ws.Range("D39").Value = "'2.802"
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'6.362"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'75.97"
$ws.Range("E41").Value = "  +6.46%  "
$ws.Range("D42").Value = "'2.018"
$ws.Range("E42").Value = "  +1.22%  "
$ws.Range("D43").Value = "'0.4261"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "'0.9989"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'0.8342"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").Value = "'101.17"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("D47").Value = "'9.444"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("D48").Value = "'7.043"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "'35.63"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").Value = "'917.40"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'0.05751"
$ws.Range("E51").Value = "  +1.89%  "
